$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rebuild the sheet as a two-column "term -> canonical vaccine name" lookup
# table. Column A keeps the original vaccine-name / alias rows (now with a
# new header), column B is new and maps every row to its canonical vaccine.
# ---------------------------------------------------------------------------

$rows = @(
    @{ A = "vaccine_key_terms";                                              ABold = $true;  ACat = $false; B = "vaccine";                     BBold = $true },
    @{ A = "Moderna mRNA-1273";                                              ACat = $true;   B = "Moderna" },
    @{ A = "Moderna";                                                        ACat = $false;  B = "Moderna" },
    @{ A = "Pfizer/BioNTech BNT162b2";                                       ACat = $true;   B = "Pfizer-BioNTech" },
    @{ A = "Pfizer";                                                         ACat = $false;  B = "Pfizer-BioNTech" },
    @{ A = "Comirnaty";                                                      ACat = $false;  B = "Pfizer-BioNTech" },
    @{ A = "Johnson & Johnson Janssen Ad26.COV2.S";                          ACat = $true;   B = "Janssen (Johnson & Johnson)" },
    @{ A = "Johnson & Johnson";                                              ACat = $false;  B = "Janssen (Johnson & Johnson)" },
    @{ A = "Janssen";                                                        ACat = $false;  B = "Janssen (Johnson & Johnson)" },
    @{ A = "Sinovac CoronaVac";                                              ACat = $true;   B = "Sinovac CoronaVac";          BCat = $true },
    @{ A = "Sinovac";                                                        ACat = $false;  B = "Sinovac CoronaVac";          BCat = $true },
    @{ A = "CoronaVac";                                                      ACat = $false;  B = "Sinovac CoronaVac";          BCat = $true },
    @{ A = "Oxford/AstraZeneca AZD1222";                                     ACat = $true;   B = "AstraZeneca";                BCat = $true },
    @{ A = "AstraZeneca";                                                    ACat = $false;  B = "AstraZeneca";                BCat = $true },
    @{ A = "Sinopharm BBIBP-CorV Vero Cells";                                ACat = $true;   B = "Sinopharm";                  BCat = $true },
    @{ A = "Sinopharm";                                                      ACat = $false;  B = "Sinopharm";                  BCat = $true },
    @{ A = "Covishield Oxford/AstraZeneca Serum Institute of India";         ACat = $true;   B = "Covishield" },
    @{ A = "Covishield";                                                     ACat = $false;  B = "Covishield" }
)

$r = 1
foreach ($row in $rows) {
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)

    $aCell.Value = $row.A
    $bCell.Value = $row.B

    # Reset to the default "Normal" style first so leftover formatting from
    # the row's previous contents (this sheet is being fully re-laid-out)
    # never bleeds into cells that should end up unstyled.
    $aCell.Style = "Normal"
    $bCell.Style = "Normal"

    if ($row.ABold) {
        $aCell.Font.Bold = $true
    } elseif ($row.ACat) {
        $aCell.Font.Color = 0
    }

    if ($row.BBold) {
        $bCell.Font.Bold = $true
    } elseif ($row.BCat) {
        $bCell.Font.Color = 0
    }

    $r = $r + 1
}

# Column A width (stored OOXML width of 47 characters).
$ws.Columns.Item(1).ColumnWidth = 46.1666666666667

# Selection, matching the saved view state.
$ws.Range("A4").Select()

# Portrait page setup.
$ws.PageSetup.Orientation = 1
